$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()
$ws.Range("H11").Select()
Write-Host "done"
